# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. before "2022-Q2"),
#    populated with the fund-holdings detail for that quarter.
# 2. Insert a corresponding new row into the "总计" summary sheet (row 2),
#    pushing the existing quarters down by one row, and renumber the
#    running index in column A.

function Set-TextCell($range, $value) {
    # Force a numeric-looking string (fund codes, percentages, ...) to be
    # stored as text, matching the source data, then drop the now-stray
    # "@" number-format so the cell is left with no explicit style.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New sheet "2022-Q4"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Add($wb.Worksheets.Item("2022-Q2"))
$ws.Name = "2022-Q4"

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$fmtSrc = $wb.Worksheets.Item("2022-Q2")
$fmtSrc.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$rows = @(
    @("004008", "中融鑫思路灵活配置混合A", "1.02", "37.86", "1.94", "0.0198", 7),
    @("004009", "中融鑫思路灵活配置混合C", "0.70", "37.86", "1.94", "0.0136", 7),
    @("012415", "德邦上证 G60 创新综合指数增强A", "0.09", "91.89", "3.56", "0.0032", 6),
    @("001412", "德邦鑫星价值灵活配置混合A", "0.13", "35.79", "2.27", "0.0030", 4),
    @("002112", "德邦鑫星价值灵活配置混合C", "0.02", "35.79", "2.27", "0.0005", 4),
    @("012416", "德邦上证 G60 创新综合指数增强C", "0.01", "91.89", "3.56", "0.0004", 6)
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $r - 2
    Set-TextCell $ws.Range("B$r") $row[0]
    Set-TextCell $ws.Range("C$r") $row[1]
    Set-TextCell $ws.Range("D$r") $row[2]
    Set-TextCell $ws.Range("E$r") $row[3]
    Set-TextCell $ws.Range("F$r") $row[4]
    Set-TextCell $ws.Range("G$r") $row[5]
    $ws.Range("H$r").Value = $row[6]
    $r = $r + 1
}

$aFmtSrc = $wb.Worksheets.Item("2022-Q2")
$aFmtSrc.Range("A2").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. "总计" summary sheet: insert the 2022-Q4 row at the top of the data
#    (row 2), pushing every other quarter down by one row. Cells are
#    rewritten directly (bottom-up) instead of using Rows.Insert so the
#    existing per-row styles (A column) stay exactly where they were.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 11
$total.Range("D7").Value = 0.43

$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 15
$total.Range("D6").Value = 2.03

$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.05

$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 1.55

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.09

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 0.04

for ($i = 0; $i -le 5; $i++) {
    $total.Range("A" + (2 + $i)).Value = $i
}

$aSrc = $wb.Worksheets.Item("总计")
$aSrc.Range("A6").Copy()
$aSrc.Range("A7").PasteSpecial(-4122)
$aSrc.Range("A7").Value = 5
